# docs/diagrams/LogicComponentClassDiagram.pptx: Update "AddressBook" to "Scheduler"
$p = $ppt.ActivePresentation

$found = $false
for ($k = 1; $k -le $p.Slides.Count; $k++) {
    $s = $p.Slides.Item($k)
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $shape = $s.Shapes.Item($i)
        if ($shape.HasTextFrame) {
            $tf = $shape.TextFrame
            if ($tf.HasText) {
                $tr = $tf.TextRange
                for ($j = 1; $j -le $tr.Paragraphs().Count; $j++) {
                    $para = $tr.Paragraphs($j)
                    $ptext = $para.Text.TrimEnd([char]13, [char]10)
                    if ($ptext -eq "AddressBook") {
                        $para.Text = "Scheduler"
                        $found = $true
                    }
                }
            }
        }
    }
}

Write-Host "Updated AddressBook -> Scheduler: $found"
